$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.713.29"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.992.24"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.752"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000320"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "4.632.98"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "4.007.15"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E16").Value = "  +6.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D20").Value = "72.515.25"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.57%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "682.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.442"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "0.0₃0858"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.78%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.18%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "2.806.75"
$ws.Range("E51").Value = "  +9.52%  "
